$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 104.67
$ws.Range("E2").Value = 55
$ws.Range("F2").Value = 14.54
$ws.Range("K2").Value = 64.09999999999999
$ws.Range("N2").Value = 52.28493729186943

# Row 3 updates
$ws.Range("D3").Value = 21.39
$ws.Range("E3").Value = 46.8
$ws.Range("F3").Value = 6.95
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 53
$ws.Range("J3").Value = 56
$ws.Range("K3").Value = 52.9
$ws.Range("N3").Value = 52.28493729186943
